$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")

# --- Add new checklist row (row 7) describing the Drillthrough scenario ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Drillthrough"
$ws.Range("D7").Value = "1.Generate a chart with some data`r`n2.Create new report page and in DrillThrough add the fields for drillthrough.`r`n3. Right click on the chart, select the Drillthrough option from the menu. "
$ws.Range("C7").Value = "Created custom menu to drillthrough from one visual to another."
$ws.Range("E7").Value = "1. On right click of the chart and selecting the Drillthrough option from the context-menu , the report will drillthrough to the newly created report page."

# C7/D7/E7 mirror the wrap-text-only formatting used by the sibling column (D) in other rows
$ws.Range("C7:E7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 75

# --- Tidy up row 2: drop the redundant "applyFont" style left on these cells ---
# (keep D2's existing wrap-text formatting untouched)
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# --- Restore the active selection on the BVT sheet ---
$ws.Range("B2").Select()
